$d = $word.ActiveDocument

# 1. Title: " Breast Cancer Imaging" -> " Cancer Imaging"
$d.Content.Find.Execute("Breast Cancer Imaging", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Cancer Imaging", 2)

# 2. Body: "the tissue was cut" -> "the gastric cancer biopsy tissues were cut"
$d.Content.Find.Execute("tissue was", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "gastric cancer biopsy tissues were", 2)

# 3. Word keeps an automatic "_GoBack" bookmark at the location of the most
#    recent edit. Move it from its old spot (an empty paragraph after
#    "Problem statement") to right after "were", matching what Word would do
#    after this edit.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$found = $d.Content
$found.Find.Execute("tissues were", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$bmRange = $d.Range($found.End, $found.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
